$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.625.02'
$ws.Range('D2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.563.91'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.24%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '210.33'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.13%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.506'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.38%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '24.89'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +5.52%  '

$ws.Range('E9').Value = '  +0.88%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0586'
$ws.Range('D10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0895'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.16%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.787.90'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.29%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.562.77'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.15%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.658.31'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.25%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.517'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.95%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.63'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.22%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.47'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.88%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '226.93'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.28%  '

$ws.Range('E19').Value = '  -0.53%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0679'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.31%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').Style = "Normal"

$ws.Range('E22').Value = '  -0.09%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.02'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.39%  '

$ws.Range('E24').Value = '  +1.16%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '151.53'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.97%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.105'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.20%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.76'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.58%  '

$ws.Range('E28').Value = '  -0.13%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.21'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.76%  '

$ws.Range('E30').Value = '  -3.72%  '

$ws.Range('E31').Value = '  -0.16%  '

$ws.Range('E32').Value = '  +0.33%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.399.75'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.76%  '

$ws.Range('E34').Value = '  -3.00%  '

$ws.Range('E35').Value = '  -2.71%  '

$ws.Range('E36').Value = '  -1.78%  '

$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.67'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.36%  '

$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.30'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.91%  '

$ws.Range('E39').Value = '  -0.37%  '

$ws.Range('E40').Value = '  +0.49%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.516'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.53%  '

$ws.Range('E42').Value = '  -0.11%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.765'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.74%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0460'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.55%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '63.81'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.84%  '

$ws.Range('E46').Value = '  -1.96%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.699.43'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.20%  '

$ws.Range('E48').Value = '  -5.26%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '84.80'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.50%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '42.28'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +5.48%  '

$ws.Range('E51').Value = '  -0.49%  '
